$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before E ("21 a 28 anos") - shifts old E..I to F..J
# ---------------------------------------------------------------------------
$ws.Range("E:E").EntireColumn.Insert()
$ws.Range("E2").Value = "21 a 28 anos"

# ---------------------------------------------------------------------------
# 2. Fill in the new life-history content.
#    The insertion order below matches the order the strings must be added
#    to the shared-strings table.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Minha irmã nasceu"
$ws.Range("B4").Value = "Minha vó faleceu"
$ws.Range("B5").Value = "Amizades? Tinha um pessoal barra pesada no meu bairro"

$ws.Range("D3").Value = "Vários, foi a melhor época. O mais marcante foi ter escolhido viver em república durante a faculdade. Conheci minha atual namorada/futura esposa"
$ws.Range("D5").Value = "O curso que escolhi, as pessoas que decidi seguir sendo próximo. Concorrer a bolsa de estudos do ciências sem fronteiras."

$ws.Range("E3").Value = "Vários também, neste período fui para o intercâmbio, finalizei a faculdade, comecei meu primeiro emprego, empreendi em um novo ramo, abandonei tudo e recomecei"
$ws.Range("E4").Value = "Meus pais se separaram, e tive um afastamento novamente de amigos, por conta de trabalho e estudos."
$ws.Range("E5").Value = "Iniciar/abandonar empregos, empreendimentos. Perdoar algumas pessoas, me afastar de outras."

$ws.Range("C3").Value = "Mudei para uma escola melhor, com uma estrutura muito legal, além de mudar de bairro, para uma casa melhor."
$ws.Range("D4").Value = "Mais uma vez me afastei de muitos amigos, pela mudança. Me separei da minha primeira namorada."
$ws.Range("C4").Value = "Perdi meus amigos pela mudança, fiquei um pouco isolado devido ao bairro ter mt menos crianças."
$ws.Range("C5").Value = "Concorrer à bolsa escolar e ser premiado, pois meus pais precisaram disso. "

# ---------------------------------------------------------------------------
# 3. Formatting
# ---------------------------------------------------------------------------

# Title cell (A1): centred, wrapped, vertically centred
$ws.Range("A1").WrapText = $true
$ws.Range("A1").VerticalAlignment = -4108

# Header row (A2 + the nine age-range headers B2:J2): centred, wrapped,
# vertically centred
$hdr = $ws.Range("A2:J2")
$hdr.HorizontalAlignment = -4108
$hdr.WrapText = $true
$hdr.VerticalAlignment = -4108

# Filler cells on row 1 next to the merged title (B1:J1)
$ws.Range("B1:J1").HorizontalAlignment = -4108
$ws.Range("B1:J1").VerticalAlignment = -4108
$ws.Range("B1:J1").WrapText = $true

# Data cells (B3:E5): centred, wrapped, vertically centred
$data = $ws.Range("B3:E5")
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108
$data.WrapText = $true

# ---------------------------------------------------------------------------
# 4. Row heights for the newly-populated rows
# ---------------------------------------------------------------------------
$ws.Rows("3").RowHeight = 114.75
$ws.Rows("4").RowHeight = 97.5
$ws.Rows("5").RowHeight = 97.5

# ---------------------------------------------------------------------------
# 5. View state: select C5, keep the sheet tab active
# ---------------------------------------------------------------------------
$ws.Range("C5").Select() | Out-Null

Write-Output "done"
